$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated AgTests (F) and AgPosit (G) values for existing rows 334-363
$updates = @(
    @{Row=334; F=204444; G=3467},
    @{Row=335; F=129742; G=2946},
    @{Row=336; F=103976; G=3325},
    @{Row=337; F=106623; G=3009},
    @{Row=338; F=223357; G=3121},
    @{Row=339; F=652372; G=5584},
    @{Row=340; F=376462; G=3235},
    @{Row=341; F=297453; G=3679},
    @{Row=342; F=180480; G=3049},
    @{Row=343; F=132207; G=2937},
    @{Row=344; F=135785; G=2505},
    @{Row=345; F=289602; G=3298},
    @{Row=346; F=654973; G=4704},
    @{Row=347; F=334278; G=2826},
    @{Row=348; F=233830; G=3240},
    @{Row=349; F=159148; G=2753},
    @{Row=350; F=128646; G=2775},
    @{Row=351; F=149280; G=2801},
    @{Row=352; F=303442; G=3506},
    @{Row=353; F=705857; G=5162},
    @{Row=354; F=303400; G=2766},
    @{Row=355; F=220726; G=3402},
    @{Row=356; F=158684; G=2856},
    @{Row=357; F=137465; G=3007},
    @{Row=358; F=160732; G=2663},
    @{Row=359; F=317671; G=3320},
    @{Row=360; F=724881; G=4910},
    @{Row=361; F=324412; G=2548},
    @{Row=362; F=217227; G=2962},
    @{Row=363; F=177897; G=2645}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 6).Value = $u.F
    $ws.Cells.Item($u.Row, 7).Value = $u.G
}

# New row 364 with data for 2021-03-03
$ws.Cells.Item(364, 1).Value = 44258
$ws.Cells.Item(364, 1).NumberFormat = $ws.Cells.Item(363, 1).NumberFormat
$ws.Cells.Item(364, 2).Value = 317159
$ws.Cells.Item(364, 3).Value = 12624
$ws.Cells.Item(364, 4).Value = 2800
$ws.Cells.Item(364, 5).Value = 7560
$ws.Cells.Item(364, 6).Value = 140488
$ws.Cells.Item(364, 7).Value = 2564
